$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell references to their new text values, applied in a way that
# preserves the existing Text cell type (these are all inline/text strings
# in the source data, e.g. percentages and decimal price strings that must
# not be reinterpreted as numbers).
$updates = [ordered]@{
    'D2' = '303.67'
    'E2' = '3.03%'
    'D3' = '43.03'
    'E3' = '6.55%'
    'D4' = '5.036'
    'E4' = '0.45%'
    'D5' = '0.07688'
    'E5' = '4.43%'
    'D6' = '4.408'
    'E6' = '2.73%'
    'E7' = '3.88%'
    'D8' = '1.061'
    'E8' = '14.69%'
    'D10' = '0.1234'
    'E10' = '4.27%'
    'D11' = '0.1867'
    'E11' = '4.17%'
    'D12' = '0.09169'
    'E12' = '4.65%'
    'D13' = '0.04165'
    'E13' = '-1.91%'
    'E14' = '-0.75%'
    'D15' = '0.001268'
    'E15' = '-0.77%'
    'D16' = '0.005769'
    'E16' = '-3.65%'
    'E17' = '1,893.97%'
    'D18' = '3.327'
    'E18' = '-1.24%'
    'D19' = '0.3345'
    'E19' = '2.35%'
    'D20' = '8.393'
    'E20' = '6.61%'
    'D21' = '0.1406'
    'E21' = '1.76%'
    'D22' = '0.2891'
    'E22' = '3.05%'
    'D23' = '0.04142'
    'E23' = '4.98%'
    'D24' = '0.001274'
    'E24' = '0.52%'
    'D25' = '0.004490'
    'E25' = '17.30%'
    'D26' = '0.0001347'
    'E26' = '9.26%'
    'D38' = '0.02459'
    'E38' = '4.53%'
    'D39' = '0.05282'
    'E39' = '3.56%'
    'D40' = '0.006004'
    'E40' = '-2.26%'
    'D41' = '0.007659'
    'E41' = '-1.44%'
    'D42' = '0.1350'
    'E42' = '4.01%'
    'D43' = '0.007337'
    'E43' = '-0.80%'
    'D44' = '0.008343'
    'E44' = '16.14%'
    'D45' = '0.3023'
    'E45' = '2.91%'
    'D46' = '0.00006628'
    'E46' = '8.58%'
    'E47' = '-0.40%'
    'B48' = 'BOLO'
    'C48' = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
    'D48' = '0.03842'
    'E48' = '-17.53%'
    'B49' = 'CoinbaseStockToken'
    'C49' = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
    'D49' = '0.004201'
    'E49' = '-0.10%'
    'E50' = '-0.40%'
    'D51' = '0.0001994'
    'E51' = '-0.40%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
